# joint version of the data for dashboard
# Adds a new "link" column (col 125 / DU) built from the project ID (col A)
# and fixes a handful of category cells whose values were out of order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$linkCol = 125   # column DU, right after DT (124)

# Header for the new column
$ws.Cells.Item(1, $linkCol).Value = "link"

# Populate the link for every data row (2..61) from the ID in column A
$lastRow = 61
for ($r = 2; $r -le $lastRow; $r++) {
    $id = $ws.Cells.Item($r, 1).Text
    $url = "https://app.zohocreator.eu/erp.forms20/erp/#Form:Projects?recLinkID=" + $id + "&viewLinkName=Project_spending"
    $ws.Cells.Item($r, $linkCol).Value = $url
}

# A few donor/implementor category cells were re-ordered; apply the exact
# corrected values (column letters kept as comments for clarity).
$fixes = @(
    @{ Row = 18; Col = 110; Value = "Consultancy" },       # DF18
    @{ Row = 18; Col = 114; Value = "NGO" },                # DJ18
    @{ Row = 18; Col = 120; Value = "NGO" },                # DP18
    @{ Row = 18; Col = 121; Value = "Consultancy" },        # DQ18

    @{ Row = 30; Col = 107; Value = "NGO" },                # DC30
    @{ Row = 30; Col = 115; Value = "Government" },         # DK30

    @{ Row = 42; Col = 98;  Value = "Private donors" },     # CT42
    @{ Row = 42; Col = 99;  Value = "Bilateral" },          # CU42

    @{ Row = 54; Col = 97;  Value = "Bilateral" },          # CS54
    @{ Row = 54; Col = 98;  Value = "Private donors" },     # CT54

    @{ Row = 59; Col = 106; Value = "NGO" },                # DB59
    @{ Row = 59; Col = 110; Value = "Social enterprise" },  # DF59

    @{ Row = 60; Col = 97;  Value = "Multi-lateral" },      # CS60
    @{ Row = 60; Col = 98;  Value = "Bilateral" },          # CT60
    @{ Row = 60; Col = 100; Value = "Bilateral" },          # CV60
    @{ Row = 60; Col = 101; Value = "Multi-lateral" },      # CW60
    @{ Row = 60; Col = 108; Value = "UN" },                 # DD60
    @{ Row = 60; Col = 109; Value = "Government" },         # DE60
    @{ Row = 60; Col = 110; Value = "NGO" },                # DF60
    @{ Row = 60; Col = 111; Value = "NGO" },                # DG60
    @{ Row = 60; Col = 114; Value = "Government" },         # DJ60
    @{ Row = 60; Col = 115; Value = "NGO" },                # DK60
    @{ Row = 60; Col = 116; Value = "Government" }          # DL60
)

foreach ($fix in $fixes) {
    $ws.Cells.Item($fix.Row, $fix.Col).Value = $fix.Value
}
